$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-05-08 Thursday" "2025-05-09 Friday"

Replace-Text "321÷8=40, 1" "482÷6=80, 2"
Replace-Text "207÷9=23, 0" "914÷3=304, 2"
Replace-Text "332÷8=41, 4" "566÷9=62, 8"
Replace-Text "733÷6=122, 1" "500÷9=55, 5"
Replace-Text "442÷7=63, 1" "836÷9=92, 8"

Replace-Text "195÷3=65, 0" "238÷6=39, 4"
Replace-Text "277÷5=55, 2" "542÷4=135, 2"
Replace-Text "668÷4=167, 0" "275÷2=137, 1"
Replace-Text "534÷8=66, 6" "793÷9=88, 1"
Replace-Text "826÷9=91, 7" "634÷2=317, 0"

Replace-Text "373÷8=46, 5" "901÷5=180, 1"
Replace-Text "378÷4=94, 2" "295÷9=32, 7"
Replace-Text "884÷2=442, 0" "601÷5=120, 1"
Replace-Text "322÷9=35, 7" "523÷5=104, 3"
Replace-Text "391÷4=97, 3" "557÷8=69, 5"

Replace-Text "771÷8=96, 3" "975÷3=325, 0"
Replace-Text "207÷6=34, 3" "547÷2=273, 1"
Replace-Text "244÷3=81, 1" "604÷6=100, 4"
Replace-Text "164÷3=54, 2" "270÷8=33, 6"
Replace-Text "706÷4=176, 2" "728÷8=91, 0"

Replace-Text "485÷8=60, 5" "579÷3=193, 0"
Replace-Text "242÷6=40, 2" "909÷3=303, 0"
Replace-Text "503÷4=125, 3" "586÷3=195, 1"
Replace-Text "504÷2=252, 0" "306÷3=102, 0"
Replace-Text "333÷4=83, 1" "494÷5=98, 4"
